# plantilla_docentes.xlsx - "Diseño vistas - pagina en general"
# Collapses the 4 separate "paralelo_a/b/c/d" columns (H:K) into a single
# "paralelo" column (H), shifts the old "ciclo" column (L) left into I,
# drops the now-unused header style on H2:H13, and moves the active
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Move the "ciclo" column (previously L) into its new position (I),
#    re-using the very same text/shared-string so nothing new is interned.
$ws.Range("I1").Value = "ciclo"

# 2) Drop the old J:L columns (paralelo_b, paralelo_c, paralelo_d headers
#    plus their "NA" data) entirely - this also shrinks the sheet
#    dimension down to A1:I13 automatically.
$ws.Range("J1:L13").ClearContents()

# 3) Rename the remaining "paralelo_a" header (H1) to just "paralelo" -
#    it is now the sole cell referencing that shared string, so the
#    string is renamed in place instead of creating a new entry.
$ws.Range("H1").Value = "paralelo"

# 4) The H2:H13 data cells used to carry an extra (redundant) explicit
#    style index; clear that formatting so they fall back to the default
#    style, matching the rest of the column.
$ws.Range("H2:H13").ClearFormats()

# 5) Update the remembered selection/active cell for the sheet view.
$ws.Range("F16").Select() | Out-Null

$wb.Save()
